$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the time_taken (column F) timestamps on the "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:33:07.814417"
$ws.Range("F3").Value  = "2021-10-05 14:33:07.814425"
$ws.Range("F4").Value  = "2021-10-05 14:33:07.814429"
$ws.Range("F5").Value  = "2021-10-05 14:33:07.814431"
$ws.Range("F6").Value  = "2021-10-05 14:33:07.814434"
$ws.Range("F7").Value  = "2021-10-05 14:33:07.814437"
$ws.Range("F8").Value  = "2021-10-05 14:33:07.814440"
$ws.Range("F9").Value  = "2021-10-05 14:33:07.814442"
$ws.Range("F10").Value = "2021-10-05 14:33:07.814445"
$ws.Range("F11").Value = "2021-10-05 14:33:07.814448"
$ws.Range("F12").Value = "2021-10-05 14:33:07.814450"
$ws.Range("F13").Value = "2021-10-05 14:33:07.814453"
$ws.Range("F14").Value = "2021-10-05 14:33:07.814455"
$ws.Range("F15").Value = "2021-10-05 14:33:07.814458"
$ws.Range("F16").Value = "2021-10-05 14:33:07.814461"

# --- Add the new "metadata" sheet, placed after the "data" sheet ---
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Match page margins used elsewhere in the workbook (inches -> points)
$meta.PageSetup.LeftMargin = 54
$meta.PageSetup.RightMargin = 54
$meta.PageSetup.TopMargin = 72
$meta.PageSetup.BottomMargin = 72
$meta.PageSetup.HeaderMargin = 36
$meta.PageSetup.FooterMargin = 36

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row
$idxCell = $meta.Range("A2")
$idxCell.Value = 0
$idxCell.Font.Bold = $true
$idxCell.Borders.LineStyle = 1
$idxCell.HorizontalAlignment = -4108
$idxCell.VerticalAlignment = -4160

$meta.Range("B2").Value = "Alternating Hemiplegia and Hemiplegic Migraine"
$meta.Range("C2").Value = 40
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "0.51"
$meta.Range("E2").Value = "2021-09-13T08:18:25.402821Z"
$meta.Range("F2").Value = "2021-10-05 14:33:07.810527"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/40/?format=json"

$ws.Activate()
